# Daily update: append the next day's row of data to the bottom of the
# "Wins Over Time" tracking table (Sheet1, columns A:D).
#
# The sheet has a running table starting at row 2 (row 1 is the header:
# Day | Chase | Bryce | Zach). Each day a new row is appended with the
# date serial in column A and the three contestants' cumulative win
# counts in columns B:D. This adds the next row, one day after the
# current last row, re-using that row's date/number formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Find the current last used row of the table and compute the next one.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row   # xlUp = -4162
$newRow = $lastRow + 1

$lastDate = $ws.Cells.Item($lastRow, 1).Value2

# New day's values.
$ws.Cells.Item($newRow, 1).Value = $lastDate + 1
$ws.Cells.Item($newRow, 2).Value = 128
$ws.Cells.Item($newRow, 3).Value = 145
$ws.Cells.Item($newRow, 4).Value = 133

# Match the date column's existing number formatting.
$ws.Cells.Item($newRow, 1).NumberFormat = $ws.Cells.Item($lastRow, 1).NumberFormat
